$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Properties")

# Append a new property row (row 4) to the "Properties" sheet, mirroring
# the existing rows' layout: ID, Name, Address, Type, Sections, CreatedAt.
# The ID looks like a plain number, so prefix it with an apostrophe
# (Excel's "store as text" marker) to keep it a text value like the other
# ID cells in the sheet; the formatting is cleared right after so the
# apostrophe marker doesn't linger as a visible quote-prefix style.
$ws.Cells.Item(4, 1).Value = "'1756657070051"
$ws.Cells.Item(4, 2).Value = "test"
$ws.Cells.Item(4, 3).Value = "123idp"
$ws.Cells.Item(4, 4).Value = "partial"
$ws.Cells.Item(4, 5).Value = '["1 floor (200)","2 floor (201)","3 floor (202)","4 floor (203)","5 floor (204)"]'
$ws.Cells.Item(4, 6).Value = "2025-08-31T16:17:50.052Z"

$ws.Range("A4:F4").ClearFormats()
